# GDE-6732 Update With New Codes Correct One Locator
$wb = $excel.ActiveWorkbook

# Writes a numeric-looking loan code as literal text into $range, keeping the
# cell's original (non quote-prefixed) number format/style - a plain
# assignment of a digit string would otherwise be auto-coerced to a number,
# and a leading "'" (needed to stop that coercion) otherwise leaves Excel's
# quote-prefix marker baked into the cell style.
function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $savedStyle
}

# --- SERV01_LoanDrawdown: update loan codes + maturity/effective dates (2018 -> 2020/2021) ---
$wsLoan = $wb.Worksheets.Item("SERV01_LoanDrawdown")

Set-TextValue $wsLoan.Range("E2") "60000485"
$wsLoan.Range("K2").Value = "'31-Oct-2020"
$wsLoan.Range("L2").Value = "'30-Nov-2020"

Set-TextValue $wsLoan.Range("E3") "60000486"
$wsLoan.Range("K3").Value = "'19-Nov-2020"
$wsLoan.Range("L3").Value = "'19-Dec-2020"

Set-TextValue $wsLoan.Range("E4") "60000479"
$wsLoan.Range("K4").Value = "'28-Nov-2020"
$wsLoan.Range("L4").Value = "'28-Dec-2020"

$wsLoan.Range("K5").Value = "'19-Feb-2021"
$wsLoan.Range("L5").Value = "'19-Mar-2021"

# Selection left on A5 for this sheet
$null = $wsLoan.Range("A5").Select()

# --- UAT04_Fees: selection moves to A2 ---
$wsFees = $wb.Worksheets.Item("UAT04_Fees")
$null = $wsFees.Range("A2").Select()

# --- UAT04_Runbook: same loan-code updates + becomes the active/selected tab ---
$wsRunbook = $wb.Worksheets.Item("UAT04_Runbook")
Set-TextValue $wsRunbook.Range("E2") "60000485"
Set-TextValue $wsRunbook.Range("E6") "60000486"
Set-TextValue $wsRunbook.Range("E9") "60000479"

# Activate the runbook sheet last so it becomes the active tab / tabSelected sheet,
# and leave the selection at I20 as recorded in the saved view state.
$null = $wsRunbook.Activate()
$null = $wsRunbook.Range("I20").Select()
